# Weekly Fruta/Hortaliza update: insert a new week's data (two rows, one
# per "Calidad" grade - Primera/Segunda) right before the current first
# data block for this date range (old row 446), pushing all subsequent
# rows down by two. This mirrors the source diff, where rows 446-489
# (old) become rows 448-491 (new) and the newly inserted rows 446-447
# carry the new week's prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 446:447 - everything from old row 446 onward
# shifts down by two rows (old 446 -> new 448, ..., old 489 -> new 491).
$ws.Rows("446:447").Insert()

# Row 446 - Betarraga, Primera
$ws.Range("A446").Value = 8
$ws.Range("B446").Value = "Terminal La Palmera de La Serena"
$ws.Range("C446").Value = "Coquimbo"
$ws.Range("D446").Value = 45132
$ws.Range("E446").Value = 4
$ws.Range("F446").Value = 100114014
$ws.Range("G446").Value = "Betarraga"
$ws.Range("H446").Value = "Sin especificar"
$ws.Range("I446").Value = "Primera"
$ws.Range("J446").Value = 1700
$ws.Range("K446").Value = 550
$ws.Range("L446").Value = 600
$ws.Range("M446").Value = 575
$ws.Range("N446").Value = "`$/paquete 3 unidades"
$ws.Range("O446").Value = "Provincia del Elquí"
$ws.Range("P446").Value = 192
$ws.Range("Q446").Value = 3
$ws.Range("R446").Value = "Hortaliza"

# Row 447 - Betarraga, Segunda
$ws.Range("A447").Value = 8
$ws.Range("B447").Value = "Terminal La Palmera de La Serena"
$ws.Range("C447").Value = "Coquimbo"
$ws.Range("D447").Value = 45132
$ws.Range("E447").Value = 4
$ws.Range("F447").Value = 100114014
$ws.Range("G447").Value = "Betarraga"
$ws.Range("H447").Value = "Sin especificar"
$ws.Range("I447").Value = "Segunda"
$ws.Range("J447").Value = 800
$ws.Range("K447").Value = 450
$ws.Range("L447").Value = 500
$ws.Range("M447").Value = 475
$ws.Range("N447").Value = "`$/paquete 3 unidades"
$ws.Range("O447").Value = "Provincia del Elquí"
$ws.Range("P447").Value = 158
$ws.Range("Q447").Value = 3
$ws.Range("R447").Value = "Hortaliza"
